# Update the cryptos price/volume sheet with the latest scraped values.
# Price cells (column D) are forced to text via a leading apostrophe so that
# values such as "30.994.16" or "239.24" are stored as literal strings
# instead of being auto-converted by Excel into numbers/dates, matching the
# original inlineStr "t" cells in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.994.16'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '''1.912.02'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('D4').Value = '''0.9981'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '''239.24'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('D6').Value = '''0.9984'
$ws.Range('D7').Value = '''0.4908'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '''0.2968'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '''0.06788'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').Value = '''1.926.04'
$ws.Range('E10').Value = '  +1.62%  '
$ws.Range('D11').Value = '''17.12'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '''0.07284'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '''90.12'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('D14').Value = '''5.127'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = '''0.6719'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '''30.931.00'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '''0.000007973'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').Value = '''0.9987'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '''2.150.31'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').Value = '''0.9987'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '''5.104'
$ws.Range('E22').Value = '  +4.95%  '
$ws.Range('D23').Value = '''208.46'
$ws.Range('E23').Value = '  +9.51%  '
$ws.Range('D24').Value = '''6.221'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').Value = '''9.667'
$ws.Range('E25').Value = '  +2.72%  '
$ws.Range('D26').Value = '''158.18'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').Value = '''19.00'
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('D28').Value = '''1.965'
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '''1.426'
$ws.Range('E29').Value = '  +1.88%  '
$ws.Range('D30').Value = '''4.330'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').Value = '''0.09181'
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').Value = '''4.037'
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('D33').Value = '''0.05184'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = '''0.7516'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').Value = '''1.120'
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').Value = '''2.689'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('D38').Value = '''2.733'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''2.118'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''0.9285'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').Value = '''0.4496'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').Value = '''106.68'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.830'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '''1.007'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('D45').Value = '''7.788'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').Value = '''0.1377'
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('D47').Value = '''66.81'
$ws.Range('E47').Value = '  +14.66%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '''35.13'
$ws.Range('E48').Value = '  +4.52%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '''0.4078'
$ws.Range('E49').Value = '  +2.52%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05903'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''8.894'
$ws.Range('E51').Value = '  +0.49%  '
